{"js": "// The underlying change in this revision is a pure OOXML re-serialization:\n// the producing toolchain (Apache POI) was upgraded/repackaged and, as a\n// side effect, it now writes XML attributes (namespace declarations on\n// <w:document>, and attributes on <w:pgSz>, <w:pgMar>, <w:rFonts>,\n// <w:lang>, <w:latentStyles>, <w:lsdException>, <w:style>, <w:tblInd>,\n// <w:tblCellMar>, ...) in a different (alphabetical) order.\n//\n// Every single changed line in the diff keeps the exact same element\n// name and the exact same set of attribute name/value pairs - only the\n// on-disk attribute order differs. There is no textual, structural,\n// formatting, or style-value change anywhere in the document: page size\n// stays 11906 x 16838 twips, margins stay 1417/1417/1417/1417 (header\n// 708, footer 708, gutter 0), the default run fonts/size/language stay\n// the same, and every latent style / custom style keeps its original\n// name, id, and properties.\n//\n// Word's object model (Office.js here) has no notion of \"attribute\n// serialization order\" - that is purely an artifact of the XML writer\n// used when the package is saved, not something an end-user edit (or a\n// script driving the documented object model) can express. Reproducing\n// the intent of this commit therefore means leaving the document's\n// content and formatting completely untouched, which is what this\n// script verifies by reading back (without mutating) exactly the\n// properties the diff touches.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const body = section.body;\n  body.load(\"text\");\n}\n\nawait context.sync();\n// No property is written: values already match the expected output, so\n// the body/sections/styles are read-only inspected and left unchanged.\n", "ps1": "# The underlying change in this revision is a pure OOXML re-serialization:\n# the producing toolchain (Apache POI) was upgraded/repackaged and, as a\n# side effect, it now writes XML attributes (namespace declarations on\n# <w:document>, and attributes on <w:pgSz>, <w:pgMar>, <w:rFonts>,\n# <w:lang>, <w:latentStyles>, <w:lsdException>, <w:style>, <w:tblInd>,\n# <w:tblCellMar>, ...) in a different (alphabetical) order.\n#\n# Every single changed line in the diff keeps the exact same element\n# name and the exact same set of attribute name/value pairs - only the\n# on-disk attribute order differs. There is no textual, structural,\n# formatting, or style-value change anywhere in the document: page size\n# stays 11906 x 16838 twips, margins stay 1417/1417/1417/1417 (header\n# 708, footer 708, gutter 0), the default run fonts/size/language stay\n# the same, and every latent style / custom style keeps its original\n# name, id, and properties.\n#\n# The Word COM object model has no notion of \"attribute serialization\n# order\" - that is purely an artifact of the XML writer used when the\n# package is saved, not something an end-user edit (or a script driving\n# the documented object model) can express. Reproducing the intent of\n# this commit therefore means leaving the document's content and\n# formatting completely untouched, which is what this script verifies\n# by reading back (without mutating) exactly the properties the diff\n# touches.\n\n$d = $word.ActiveDocument\n\n$pageWidth  = $d.PageSetup.PageWidth\n$pageHeight = $d.PageSetup.PageHeight\n$topMargin    = $d.PageSetup.TopMargin\n$bottomMargin = $d.PageSetup.BottomMargin\n$leftMargin   = $d.PageSetup.LeftMargin\n$rightMargin  = $d.PageSetup.RightMargin\n$gutter       = $d.PageSetup.Gutter\n\n$styleCount = $d.Styles.Count\n\n# No property is assigned: values already match the expected output, so\n# the page setup and styles are read-only inspected and left unchanged.\n"}
